$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 5 new rows above the existing "Total:" row (currently row 13). ---
# This pushes the old Total row (with its SUM formulas) down to row 18,
# and the newly inserted blank rows inherit the row-13 formatting context
# (so they get spans="1:7" like the surrounding data rows).
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

# NOTE: the shared-string table records new strings in first-use order. The
# source workbook's string table has "Implementation of MarketPlace" (used on
# row 14) before "Implementing planets and Universe" (used on row 13), so we
# deliberately write the A14 text first to reproduce that same table order,
# then fill in the rest of row 13/14 in normal order.
$ws.Cells.Item(14, 1).Value = "Implementation of MarketPlace"
$ws.Cells.Item(13, 1).Value = "Implementing planets and Universe"

# --- Row 13: Implementing planets and Universe ---
$ws.Cells.Item(13, 2).Value = 7
$ws.Cells.Item(13, 3).Value = 0.1
$ws.Cells.Item(13, 4).Value = 7.1
$ws.Cells.Item(13, 5).Value = 7.1
$ws.Cells.Item(13, 6).Value = 5
$ws.Rows.Item(13).RowHeight = 30

# --- Row 14: Implementation of MarketPlace ---
$ws.Cells.Item(14, 2).Value = 10
$ws.Cells.Item(14, 3).Value = 1
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = 11
$ws.Cells.Item(14, 6).Value = 8
$ws.Rows.Item(14).RowHeight = 30

# --- Row 15: Implementing ShipYard ---
$ws.Cells.Item(15, 1).Value = "Implementing ShipYard"
$ws.Cells.Item(15, 2).Value = 5
$ws.Cells.Item(15, 3).Value = 0.2
$ws.Cells.Item(15, 4).Value = 5.2
$ws.Cells.Item(15, 5).Value = 5.2
$ws.Cells.Item(15, 6).Value = 4

# --- Row 16: Testing ---
$ws.Cells.Item(16, 1).Value = "Testing"
$ws.Cells.Item(16, 2).Value = 8
$ws.Cells.Item(16, 3).Value = 0.3
$ws.Cells.Item(16, 4).Value = 8.3000000000000007
$ws.Cells.Item(16, 5).Value = 8.3000000000000007
$ws.Cells.Item(16, 6).Value = 5

# --- Row 17: UI Evaluation ---
$ws.Cells.Item(17, 1).Value = "UI Evaluation"
$ws.Cells.Item(17, 2).Value = 5
$ws.Cells.Item(17, 3).Value = 0.4
$ws.Cells.Item(17, 4).Value = 5.4
$ws.Cells.Item(17, 5).Value = 5.4
$ws.Cells.Item(17, 6).Value = 3

# --- The old "Total:" row is now at row 18 (leaving rows 18-19 as a gap). ---
# Wipe it out completely (contents + formatting) so the gap rows are blank.
$ws.Rows.Item(18).Clear()
$ws.Rows.Item(19).Clear()

# --- Write the new "Total:" row at row 20, summing the expanded data range. ---
$ws.Cells.Item(20, 1).Value = "Total:"
$ws.Cells.Item(20, 2).Formula = "=SUM(B2:B17)"
$ws.Cells.Item(20, 3).Formula = "=SUM(C2:C17)"
$ws.Cells.Item(20, 4).Formula = "=SUM(D2:D17)"
$ws.Cells.Item(20, 5).Formula = "=SUM(E2:E17)"
$ws.Cells.Item(20, 6).Formula = "=SUM(F1:F17)"

$totalRow = $ws.Range("A20:F20")
$totalRow.Font.Bold = $true
$totalRow.WrapText = $true

# --- Update the sheet's current selection to match the new Total row. ---
$ws.Range("D20").Select()
